# A new order (Order #7) came in at 2026-01-13 18:48 and was logged at the
# top of the "Daily Orders" log (row 2), pushing all earlier rows down by
# one. The "Summary" sheet's aggregate counters are updated to reflect the
# new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily Orders")

# Insert a brand-new row above the current row 2 (the previous top order,
# order #6) so that everything currently in rows 2-7 shifts down to rows
# 3-8, and row 2 becomes free for the new order's data.
$ws.Rows.Item(2).Insert()

# Fill in the data for new order #7 in the freshly inserted row 2.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "2026-01-13 18:48"
$ws.Range("C2").Value = "Sagar Borse"
# Phone numbers are stored as text (like the rest of the column), so force
# a text number format before writing the digit string.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "7588930329"
$ws.Range("E2").Value = "Test2,"
$ws.Range("F2").Value = "Kite Haldi Kunku Set x10"
$ws.Range("G2").Value = 300
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"

# Update the "Summary" sheet's aggregate counters for the new order.
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A2").Value = 7     # Total Orders: 6 -> 7
$ws2.Range("B2").Value = 6     # New: 5 -> 6
$ws2.Range("G2").Value = 325   # Total Revenue: 25 -> 325
